# Update countries & provincias Spain
# Applies the 22:01 data refresh on top of the 20:44 snapshot:
#   - Guinea overtakes Malasia (rows 98/99 swap rank)
#   - Islas Malvinas overtakes Montserrat (rows 214/215 swap rank)
#   - Refreshed totals for Estados Unidos (row 4), Irlanda (row 73),
#     Republica del Chad (row 166) and Eritrea (row 181)
#   - Timestamp caption updated to 22:01

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp caption (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Septiembre de 2020 a las 22:01"

# --- Row 4: Estados Unidos ---
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 7067804
$ws.Range("C4").Value = 21588
$ws.Range("D4").Value = 4320753
$ws.Range("E4").Value = 2541936
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 613
$ws.Range("H4").Value = 205115

# --- Row 73: Irlanda ---
$ws.Range("A73").Value = "Irlanda"
$ws.Range("B73").Value = 33444
$ws.Range("C73").Value = 323
$ws.Range("D73").Value = 23364
$ws.Range("E73").Value = 8288
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 1792

# --- Row 98: Guinea now outranks Malasia ---
$ws.Range("A98").Value = "Guinea"
$ws.Range("B98").Value = 10387
$ws.Range("C98").Value = 43
$ws.Range("D98").Value = 9780
$ws.Range("E98").Value = 542
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 65

# --- Row 99: Malasia drops below Guinea ---
$ws.Range("A99").Value = "Malasia"
$ws.Range("B99").Value = 10358
$ws.Range("C99").Value = 82
$ws.Range("D99").Value = 9563
$ws.Range("E99").Value = 665
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 130

# --- Row 166: Republica del Chad ---
$ws.Range("A166").Value = "Republica del Chad"
$ws.Range("B166").Value = 1155
$ws.Range("C166").Value = 2
$ws.Range("D166").Value = 967
$ws.Range("E166").Value = 107
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 81

# --- Row 181: Eritrea ---
$ws.Range("A181").Value = "Eritrea"
$ws.Range("B181").Value = 364
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 309
$ws.Range("E181").Value = 55
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0

# --- Row 214: Islas Malvinas now outranks Montserrat ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# --- Row 215: Montserrat drops below Islas Malvinas ---
$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
